# Ajustes de Inventario Ok
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Highlight D9 (Inventario Físico Paso 2) in yellow, like the other
#    highlighted step cells (D2, D4, D6, D8 ...)
$ws.Range("D9").Interior.Color = 65535

# 2) Add a note in E9 about missing print functionality for the listing
$ws.Range("E9").Value = "Falta poder imprimir el listado"

# 3) Insert a new row of content "Inventarios con filtros" before the
#    existing D13 ("Stock"), pushing D13:D24 down to D14:D25.
#    First, extend formatting from D13 down onto the newly used D25 cell.
$ws.Range("D13").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# Shift values down one row at a time, starting from the bottom so
# that source values are not overwritten before they are read.
$ws.Range("D25").Value = $ws.Range("D24").Value()
$ws.Range("D24").Value = $ws.Range("D23").Value()
$ws.Range("D23").Value = $ws.Range("D22").Value()
$ws.Range("D22").Value = $ws.Range("D21").Value()
$ws.Range("D21").Value = $ws.Range("D20").Value()
$ws.Range("D20").Value = $ws.Range("D19").Value()
$ws.Range("D19").Value = $ws.Range("D18").Value()
$ws.Range("D18").Value = $ws.Range("D17").Value()
$ws.Range("D17").Value = $ws.Range("D16").Value()
$ws.Range("D16").Value = $ws.Range("D15").Value()
$ws.Range("D15").Value = $ws.Range("D14").Value()
$ws.Range("D14").Value = $ws.Range("D13").Value()
$ws.Range("D13").Value = "Inventarios con filtros"

# 4) Update the active selection to D14 (reflecting the edited area)
$ws.Range("D14").Select()
